$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DLC_List")

# Make DLC_List the active tab (config sheet loses tabSelected automatically)
$ws.Activate()

# Bump the seed count per wind speed (8 seeds instead of 6): update the
# "[xxx:1:xxx6]" -> "[xxx:1:xxx8]" labels in column J, rows 2-12.
# These are written first so the new shared strings land before "512".
$jvals = @("[401:1:408]", "[601:1:608]", "[801:1:808]", "[1001:1:1008]", "[1201:1:1208]", "[1401:1:1408]", "[1601:1:1608]", "[1801:1:1808]", "[2001:1:2008]", "[2201:1:2208]", "[2401:1:2408]")
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 10).Value = $jvals[$r - 2]
}

# Duration column (H) changes from 600 to 512 for every DLC row (dbg for Rosco)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 8).Value = "512"
}

# Leave the selection on H12 as the last interacted cell
$ws.Range("H12").Select()
